# Author's intent (see commit message "stress list, appender, layout"):
#  1. Merge the "stress-tes" / "t " runs (split apart by a stray
#     _GoBack bookmark) back into a single "stress-test " run, dropping
#     the now-redundant bookmark.
#  2. Recolor the appender-comparison list (MemAppender/LinkedList/
#     ArrayList/ConsoleAppender/FileAppender) green, matching the
#     existing "different combinations" / "stress-test" styling used
#     elsewhere in the document.
#  3. Recolor the "PatternLayout and VelocityLayout" phrase the same
#     green.

$d = $word.ActiveDocument
$green = 5287936  # RGB(0x00, 0xB0, 0x50) -> 0x00B050 as a VBA/COM color long

# 1. "stress-tes" + hidden _GoBack bookmark + "t " -> single run "stress-test "
$r = $d.Content
$r.Find.Execute("stress-test ", $false, $false, $false, $false, $false, $true, 1, $false, "stress-test ", 2) | Out-Null

# 2. Colour the appender comparison list green (leaves the trailing
#    " - measure time..." text uncoloured, matching the diff).
$r = $d.Content
$r.Find.Execute("MemAppender using a LinkedList, MemAppender using an ArrayList, ConsoleAppender and FileAppender ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Color = $green

# 3. Colour "PatternLayout and VelocityLayout" green.
$r = $d.Content
$r.Find.Execute("PatternLayout and VelocityLayout", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Color = $green
